$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet): update F3:F6 attendee counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 275
$ws1.Range("F4").Value = 2705
$ws1.Range("F5").Value = 52
$ws1.Range("F6").Value = 578

# Sheet "全部类型" (fourth sheet): update F5:F8 attendee counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 275
$ws4.Range("F6").Value = 2705
$ws4.Range("F7").Value = 52
$ws4.Range("F8").Value = 578
